# Fruta / hortaliza, semanal
#
# A new weekly price report (2022-02-23, serial 44615) is inserted for
# Comercializadora del Agro de Limarí - Frutilla, pushing all the existing
# weekly blocks (rows 320-340) down by 3 rows. The oldest block that falls
# off the bottom (2020-11-25, serial 44160) lands naturally at the new
# rows 341-343 once the insert shifts everything down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the first data row of this block (row 320),
# shifting rows 320:340 down to 323:343 and carrying their values/styles.
$ws.Rows("320:322").Insert()

# --- New row 320: Especial, week of 2022-02-23 ---
$ws.Cells.Item(320,1).Value = 2
$ws.Cells.Item(320,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(320,3).Value = "Coquimbo"
$ws.Cells.Item(320,4).Value = 44615
$ws.Cells.Item(320,5).Value = 4
$ws.Cells.Item(320,6).Value = "Fruta"
$ws.Cells.Item(320,7).Value = 100101
$ws.Cells.Item(320,8).Value = "Berries"
$ws.Cells.Item(320,9).Value = 100112025
$ws.Cells.Item(320,10).Value = "Frutilla"
$ws.Cells.Item(320,11).Value = "Sin especificar"
$ws.Cells.Item(320,12).Value = "Especial"
$ws.Cells.Item(320,13).Value = 500
$ws.Cells.Item(320,14).Value = 12500
$ws.Cells.Item(320,15).Value = 13000
$ws.Cells.Item(320,16).Value = 12750
$ws.Cells.Item(320,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(320,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(320,19).Value = 1821
$ws.Cells.Item(320,20).Value = 7

# --- New row 321: Primera, week of 2022-02-23 ---
$ws.Cells.Item(321,1).Value = 2
$ws.Cells.Item(321,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(321,3).Value = "Coquimbo"
$ws.Cells.Item(321,4).Value = 44615
$ws.Cells.Item(321,5).Value = 4
$ws.Cells.Item(321,6).Value = "Fruta"
$ws.Cells.Item(321,7).Value = 100101
$ws.Cells.Item(321,8).Value = "Berries"
$ws.Cells.Item(321,9).Value = 100112025
$ws.Cells.Item(321,10).Value = "Frutilla"
$ws.Cells.Item(321,11).Value = "Sin especificar"
$ws.Cells.Item(321,12).Value = "Primera"
$ws.Cells.Item(321,13).Value = 500
$ws.Cells.Item(321,14).Value = 10500
$ws.Cells.Item(321,15).Value = 11000
$ws.Cells.Item(321,16).Value = 10750
$ws.Cells.Item(321,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(321,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(321,19).Value = 1536
$ws.Cells.Item(321,20).Value = 7

# --- New row 322: Segunda, week of 2022-02-23 ---
$ws.Cells.Item(322,1).Value = 2
$ws.Cells.Item(322,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(322,3).Value = "Coquimbo"
$ws.Cells.Item(322,4).Value = 44615
$ws.Cells.Item(322,5).Value = 4
$ws.Cells.Item(322,6).Value = "Fruta"
$ws.Cells.Item(322,7).Value = 100101
$ws.Cells.Item(322,8).Value = "Berries"
$ws.Cells.Item(322,9).Value = 100112025
$ws.Cells.Item(322,10).Value = "Frutilla"
$ws.Cells.Item(322,11).Value = "Sin especificar"
$ws.Cells.Item(322,12).Value = "Segunda"
$ws.Cells.Item(322,13).Value = 400
$ws.Cells.Item(322,14).Value = 8500
$ws.Cells.Item(322,15).Value = 9000
$ws.Cells.Item(322,16).Value = 8750
$ws.Cells.Item(322,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(322,18).Value = "Provincia de Melipilla"
$ws.Cells.Item(322,19).Value = 1250
$ws.Cells.Item(322,20).Value = 7
